$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the label in A2 to include "(mg)"
$ws.Range("A2").Value = "15N-NH4NO3 (mg)"

# 2) Update the formula in E6 to divide by 100 as well (scale correction)
$ws.Range("E6").Formula = "=C6/(20.99*100)"

# 3) Move the active selection to E7 (cosmetic, matches author's final cursor position)
$ws.Range("E7").Select()
